$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, pushing the existing rows 23-50 down to 24-51.
$ws.Rows("23:23").Insert()

# Fill the newly inserted row 23 with the new record's data.
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C23").Value = "Arica y Parinacota"
$ws.Range("D23").Value = 44413
$ws.Range("E23").Value = 15
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100102
$ws.Range("H23").Value = "Cítricos"
$ws.Range("I23").Value = 100102005
$ws.Range("J23").Value = "Naranja"
$ws.Range("K23").Value = "Lane Late"
$ws.Range("L23").Value = "Segunda"
$ws.Range("M23").Value = 270
$ws.Range("N23").Value = 550
$ws.Range("O23").Value = 600
$ws.Range("P23").Value = 575
$ws.Range("Q23").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R23").Value = "Región Metropolitana"
$ws.Range("S23").Value = 575
$ws.Range("T23").Value = 1
